$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" = strikeouts) values for rows 2-9 per regenerated save data
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 8
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 9
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 5
$ws.Range("G9").Value = 4
